# Apply updated cryptocurrency price/volume figures to columns D (Price) and E (Volume(1h)).
# A leading apostrophe forces Excel to treat the assigned text as a literal string rather
# than auto-converting number-like values (e.g. "572.53") to floating point numbers, which
# would otherwise introduce binary rounding artifacts. Resetting the range Style back to
# "Normal" afterwards keeps the cell formatting identical to the original (unstyled) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '''63.579.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Formula = '''2.600.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Formula = '''  -1.67%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Formula = '''  -0.01%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Formula = '''572.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Formula = '''  -3.92%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Formula = '''154.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Formula = '''  -1.17%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D8').Formula = '''0.628'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Formula = '''  -0.33%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Formula = '''  -5.35%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Formula = '''5.76'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Formula = '''  -0.82%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Formula = '''0.382'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Formula = '''  -3.11%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Formula = '''  -0.56%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Formula = '''28.07'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Formula = '''  -1.94%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Formula = '''3.072.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Formula = '''  -1.38%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Formula = '''0.0000182'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Formula = '''  -7.35%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Formula = '''63.436.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Formula = '''  -2.91%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Formula = '''2.608.11'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Formula = '''  -1.65%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Formula = '''12.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Formula = '''  -4.10%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Formula = '''  -2.33%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Formula = '''7.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Formula = '''  +0.61%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Formula = '''341.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Formula = '''  -1.82%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Formula = '''  -0.17%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Formula = '''  -2.81%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Formula = '''1.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Formula = '''  +2.47%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Formula = '''  -4.59%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Formula = '''9.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Formula = '''  -3.94%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Formula = '''575.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Formula = '''  +9.19%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Formula = '''1.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Formula = '''  -0.31%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Formula = '''  +2.02%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Formula = '''  -2.52%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Formula = '''7.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Formula = '''  -0.24%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Formula = '''2.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Formula = '''  -2.44%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Formula = '''  -3.04%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Formula = '''6.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Formula = '''  +0.73%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Formula = '''5.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Formula = '''  -1.61%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Formula = '''  -2.77%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Formula = '''19.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Formula = '''  -2.39%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Formula = '''  +0.01%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Formula = '''151.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Formula = '''  -2.19%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Formula = '''1.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Formula = '''  -3.54%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Formula = '''  -0.04%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Formula = '''41.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Formula = '''  -2.54%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Formula = '''155.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Formula = '''  -2.98%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Formula = '''2.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Formula = '''  +3.49%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Formula = '''  -2.85%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Formula = '''22.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Formula = '''  +1.24%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Formula = '''  -2.19%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Formula = '''0.101'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Formula = '''  +1.77%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Formula = '''0.625'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Formula = '''  -1.22%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Formula = '''  -1.82%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Formula = '''  -3.86%  '
$ws.Range('E51').Style = 'Normal'
